$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.822.15'
$ws.Range('E2').Value = '  +0.92%  '

$ws.Range('D3').Value = '3.128.15'
$ws.Range('E3').Value = '  +1.00%  '

$ws.Range('E4').Value = '  +0.03%  '

$r = $ws.Range('D5')
$origStyle = $r.Style
$r.Value = "'533.34"
$r.Style = $origStyle
$ws.Range('E5').Value = '  +1.73%  '

$r = $ws.Range('D6')
$origStyle = $r.Style
$r.Value = "'138.89"
$r.Style = $origStyle
$ws.Range('E6').Value = '  +1.68%  '

$ws.Range('E7').Value = '  -0.09%  '

$ws.Range('D8').Value = '3.124.95'
$ws.Range('E8').Value = '  +1.04%  '

$ws.Range('E9').Value = '  +6.13%  '

$ws.Range('E10').Value = '  +0.30%  '

$ws.Range('E11').Value = '  +0.77%  '

$ws.Range('E12').Value = '  +4.41%  '

$ws.Range('D13').Value = '3.667.37'
$ws.Range('E13').Value = '  +1.04%  '

$ws.Range('E14').Value = '  +1.91%  '

$r = $ws.Range('D15')
$origStyle = $r.Style
$r.Value = "'25.87"
$r.Style = $origStyle
$ws.Range('E15').Value = '  +2.55%  '

$ws.Range('E16').Value = '  +1.12%  '

$ws.Range('D17').Value = '57.926.66'
$ws.Range('E17').Value = '  +0.97%  '

$ws.Range('D18').Value = '3.128.88'
$ws.Range('E18').Value = '  +1.32%  '

$ws.Range('E19').Value = '  +2.85%  '

$r = $ws.Range('D20')
$origStyle = $r.Style
$r.Value = "'12.70"
$r.Style = $origStyle
$ws.Range('E20').Value = '  +2.88%  '

$r = $ws.Range('D21')
$origStyle = $r.Style
$r.Value = "'8.09"
$r.Style = $origStyle
$ws.Range('E21').Value = '  +3.20%  '

$r = $ws.Range('D22')
$origStyle = $r.Style
$r.Value = "'367.61"
$r.Style = $origStyle
$ws.Range('E22').Value = '  +5.91%  '

$ws.Range('E23').Value = '  -0.03%  '

$ws.Range('E24').Value = '  -1.86%  '

$ws.Range('E25').Value = '  +2.50%  '

$ws.Range('E26').Value = '  +1.63%  '

$ws.Range('E27').Value = '  +1.24%  '

$ws.Range('E28').Value = '  +0.12%  '

$ws.Range('D29').Value = '0.0₃0864'
$ws.Range('E29').Value = '  -2.52%  '

$r = $ws.Range('D30')
$origStyle = $r.Style
$r.Value = "'7.32"
$r.Style = $origStyle
$ws.Range('E30').Value = '  -0.25%  '

$ws.Range('E31').Value = '  +0.56%  '

$ws.Range('E32').Value = '  +1.68%  '

$r = $ws.Range('D33')
$origStyle = $r.Style
$r.Value = "'21.45"
$r.Style = $origStyle
$ws.Range('E33').Value = '  +3.46%  '

$ws.Range('E34').Value = '  +5.25%  '

$ws.Range('E35').Value = '  +3.18%  '

$r = $ws.Range('D36')
$origStyle = $r.Style
$r.Value = "'159.43"
$r.Style = $origStyle
$ws.Range('E36').Value = '  +0.29%  '

$ws.Range('E37').Value = '  +0.75%  '

$r = $ws.Range('D38')
$origStyle = $r.Style
$r.Value = "'1.29"
$r.Style = $origStyle
$ws.Range('E38').Value = '  +5.61%  '

$r = $ws.Range('D39')
$origStyle = $r.Style
$r.Value = "'25.44"
$r.Style = $origStyle
$ws.Range('E39').Value = '  -1.07%  '

$ws.Range('E40').Value = '  +4.78%  '

$r = $ws.Range('D41')
$origStyle = $r.Style
$r.Value = "'0.0672"
$r.Style = $origStyle
$ws.Range('E41').Value = '  +2.32%  '

$ws.Range('D42').Value = '2.528.57'
$ws.Range('E42').Value = '  +6.65%  '

$r = $ws.Range('D43')
$origStyle = $r.Style
$r.Value = "'4.10"
$r.Style = $origStyle
$ws.Range('E43').Value = '  +0.25%  '

$ws.Range('E44').Value = '  +0.54%  '

$r = $ws.Range('D45')
$origStyle = $r.Style
$r.Value = "'37.79"
$r.Style = $origStyle
$ws.Range('E45').Value = '  +3.52%  '

$ws.Range('E46').Value = '  +1.78%  '

$ws.Range('E47').Value = '  -0.03%  '

$r = $ws.Range('D48')
$origStyle = $r.Style
$r.Value = "'0.980"
$r.Style = $origStyle
$ws.Range('E48').Value = '  +1.32%  '

$ws.Range('E49').Value = '  +3.24%  '

$r = $ws.Range('D50')
$origStyle = $r.Style
$r.Value = "'19.76"
$r.Style = $origStyle
$ws.Range('E50').Value = '  +0.65%  '

$r = $ws.Range('D51')
$origStyle = $r.Style
$r.Value = "'0.746"
$r.Style = $origStyle
$ws.Range('E51').Value = '  -0.98%  '
